# Apply updated odds values for Row 7 and Row 8 on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 updates
$ws.Range("M7").Value = 1.13
$ws.Range("N7").Value = 6

# Row 8 updates
$ws.Range("G8").Value = 1.5
$ws.Range("H8").Value = 4.2
$ws.Range("I8").Value = 6.25
$ws.Range("J8").Value = 2.05
$ws.Range("K8").Value = 2.3
$ws.Range("L8").Value = 6
$ws.Range("N8").Value = 13
$ws.Range("Q8").Value = 1.83
$ws.Range("R8").Value = 2.03
$ws.Range("U8").Value = 1.91
$ws.Range("V8").Value = 1.8
$ws.Range("X8").Value = 7
$ws.Range("Y8").Value = 8.5
$ws.Range("Z8").Value = 10
$ws.Range("AA8").Value = 13
$ws.Range("AB8").Value = 26
$ws.Range("AC8").Value = 11
$ws.Range("AD8").Value = 8
$ws.Range("AF8").Value = 51
$ws.Range("AI8").Value = 19
$ws.Range("AJ8").Value = 67
$ws.Range("AK8").Value = 41
$ws.Range("AM8").Value = 301
$ws.Range("AO8").Value = 7.5
$ws.Range("AQ8").Value = 23
$ws.Range("AW8").Value = 7.5
$ws.Range("AZ8").Value = 126

$wb.Save()
